$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-28 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-03-01 Sunday", 2) | Out-Null
$d.Content.Find.Execute("24+32=", $true, $false, $false, $false, $false, $true, 1, $false, "40+21=", 2) | Out-Null
$d.Content.Find.Execute("40+47=", $true, $false, $false, $false, $false, $true, 1, $false, "24+70=", 2) | Out-Null
$d.Content.Find.Execute("5+25=", $true, $false, $false, $false, $false, $true, 1, $false, "48-20=", 2) | Out-Null
$d.Content.Find.Execute("53-47=", $true, $false, $false, $false, $false, $true, 1, $false, "64+7=", 2) | Out-Null
$d.Content.Find.Execute("66+20=", $true, $false, $false, $false, $false, $true, 1, $false, "12+7=", 2) | Out-Null
$d.Content.Find.Execute("87+9=", $true, $false, $false, $false, $false, $true, 1, $false, "64-25=", 2) | Out-Null
$d.Content.Find.Execute("88+1=", $true, $false, $false, $false, $false, $true, 1, $false, "98-38=", 2) | Out-Null
$d.Content.Find.Execute("5+35=", $true, $false, $false, $false, $false, $true, 1, $false, "45+6=", 2) | Out-Null
$d.Content.Find.Execute("33-1=", $true, $false, $false, $false, $false, $true, 1, $false, "87-9=", 2) | Out-Null
$d.Content.Find.Execute("78-15=", $true, $false, $false, $false, $false, $true, 1, $false, "36+49=", 2) | Out-Null
$d.Content.Find.Execute("81+16=", $true, $false, $false, $false, $false, $true, 1, $false, "1+31=", 2) | Out-Null
$d.Content.Find.Execute("74-67=", $true, $false, $false, $false, $false, $true, 1, $false, "92-77=", 2) | Out-Null
$d.Content.Find.Execute("46-38=", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=", 2) | Out-Null
$d.Content.Find.Execute("18+24=", $true, $false, $false, $false, $false, $true, 1, $false, "44+34=", 2) | Out-Null
$d.Content.Find.Execute("38+44=", $true, $false, $false, $false, $false, $true, 1, $false, "67+26=", 2) | Out-Null
$d.Content.Find.Execute("51+20=", $true, $false, $false, $false, $false, $true, 1, $false, "18+42=", 2) | Out-Null
$d.Content.Find.Execute("49-44=", $true, $false, $false, $false, $false, $true, 1, $false, "87+0=", 2) | Out-Null
$d.Content.Find.Execute("40+18=", $true, $false, $false, $false, $false, $true, 1, $false, "15+11=", 2) | Out-Null
$d.Content.Find.Execute("90-60=", $true, $false, $false, $false, $false, $true, 1, $false, "13+2=", 2) | Out-Null
$d.Content.Find.Execute("80-11=", $true, $false, $false, $false, $false, $true, 1, $false, "22-0=", 2) | Out-Null
$d.Content.Find.Execute("28-14=", $true, $false, $false, $false, $false, $true, 1, $false, "41-28=", 2) | Out-Null
$d.Content.Find.Execute("93-25=", $true, $false, $false, $false, $false, $true, 1, $false, "45+27=", 2) | Out-Null
$d.Content.Find.Execute("96-45=", $true, $false, $false, $false, $false, $true, 1, $false, "85-55=", 2) | Out-Null
$d.Content.Find.Execute("55+8=", $true, $false, $false, $false, $false, $true, 1, $false, "23+13=", 2) | Out-Null
$d.Content.Find.Execute("38+50=", $true, $false, $false, $false, $false, $true, 1, $false, "64+14=", 2) | Out-Null
$d.Content.Find.Execute("76-42=", $true, $false, $false, $false, $false, $true, 1, $false, "16-16=", 2) | Out-Null
$d.Content.Find.Execute("7+13=", $true, $false, $false, $false, $false, $true, 1, $false, "35-18=", 2) | Out-Null
$d.Content.Find.Execute("32-28=", $true, $false, $false, $false, $false, $true, 1, $false, "67-50=", 2) | Out-Null
$d.Content.Find.Execute("12+13=", $true, $false, $false, $false, $false, $true, 1, $false, "27+68=", 2) | Out-Null
$d.Content.Find.Execute("39+46=", $true, $false, $false, $false, $false, $true, 1, $false, "82-40=", 2) | Out-Null
$d.Content.Find.Execute("37+53=", $true, $false, $false, $false, $false, $true, 1, $false, "92-86=", 2) | Out-Null
$d.Content.Find.Execute("15+83=", $true, $false, $false, $false, $false, $true, 1, $false, "72-45=", 2) | Out-Null
$d.Content.Find.Execute("11+29=", $true, $false, $false, $false, $false, $true, 1, $false, "32+25=", 2) | Out-Null
$d.Content.Find.Execute("76+17=", $true, $false, $false, $false, $false, $true, 1, $false, "99-50=", 2) | Out-Null
$d.Content.Find.Execute("82+8=", $true, $false, $false, $false, $false, $true, 1, $false, "46+26=", 2) | Out-Null
$d.Content.Find.Execute("58-24=", $true, $false, $false, $false, $false, $true, 1, $false, "48+17=", 2) | Out-Null
$d.Content.Find.Execute("34+13=", $true, $false, $false, $false, $false, $true, 1, $false, "22+40=", 2) | Out-Null
$d.Content.Find.Execute("76-44=", $true, $false, $false, $false, $false, $true, 1, $false, "84-42=", 2) | Out-Null
$d.Content.Find.Execute("72-41=", $true, $false, $false, $false, $false, $true, 1, $false, "55-9=", 2) | Out-Null
$d.Content.Find.Execute("90-80=", $true, $false, $false, $false, $false, $true, 1, $false, "41-0=", 2) | Out-Null
$d.Content.Find.Execute("61-20=", $true, $false, $false, $false, $false, $true, 1, $false, "44+26=", 2) | Out-Null
$d.Content.Find.Execute("24+63=", $true, $false, $false, $false, $false, $true, 1, $false, "44-10=", 2) | Out-Null
$d.Content.Find.Execute("4+47=", $true, $false, $false, $false, $false, $true, 1, $false, "78+3=", 2) | Out-Null
$d.Content.Find.Execute("11+15=", $true, $false, $false, $false, $false, $true, 1, $false, "62-41=", 2) | Out-Null
$d.Content.Find.Execute("74+14=", $true, $false, $false, $false, $false, $true, 1, $false, "80-9=", 2) | Out-Null
$d.Content.Find.Execute("47-12=", $true, $false, $false, $false, $false, $true, 1, $false, "10+71=", 2) | Out-Null
$d.Content.Find.Execute("23-2=", $true, $false, $false, $false, $false, $true, 1, $false, "49+9=", 2) | Out-Null
$d.Content.Find.Execute("52-48=", $true, $false, $false, $false, $false, $true, 1, $false, "0+52=", 2) | Out-Null
$d.Content.Find.Execute("5+28=", $true, $false, $false, $false, $false, $true, 1, $false, "97-91=", 2) | Out-Null
$d.Content.Find.Execute("93-87=", $true, $false, $false, $false, $false, $true, 1, $false, "45+41=", 2) | Out-Null
$d.Content.Find.Execute("11+31=", $true, $false, $false, $false, $false, $true, 1, $false, "67+2=", 2) | Out-Null
$d.Content.Find.Execute("93-34=", $true, $false, $false, $false, $false, $true, 1, $false, "64-12=", 2) | Out-Null
$d.Content.Find.Execute("6+63=", $true, $false, $false, $false, $false, $true, 1, $false, "92-78=", 2) | Out-Null
$d.Content.Find.Execute("11+69=", $true, $false, $false, $false, $false, $true, 1, $false, "48-34=", 2) | Out-Null
$d.Content.Find.Execute("1+14=", $true, $false, $false, $false, $false, $true, 1, $false, "36+51=", 2) | Out-Null
$d.Content.Find.Execute("14+25=", $true, $false, $false, $false, $false, $true, 1, $false, "67-6=", 2) | Out-Null
$d.Content.Find.Execute("17+79=", $true, $false, $false, $false, $false, $true, 1, $false, "61-35=", 2) | Out-Null
$d.Content.Find.Execute("39-20=", $true, $false, $false, $false, $false, $true, 1, $false, "42+19=", 2) | Out-Null
$d.Content.Find.Execute("73+6=", $true, $false, $false, $false, $false, $true, 1, $false, "80+2=", 2) | Out-Null
$d.Content.Find.Execute("32-13=", $true, $false, $false, $false, $false, $true, 1, $false, "99-35=", 2) | Out-Null
$d.Content.Find.Execute("40-9=", $true, $false, $false, $false, $false, $true, 1, $false, "40+11=", 2) | Out-Null
$d.Content.Find.Execute("30+40=", $true, $false, $false, $false, $false, $true, 1, $false, "52-22=", 2) | Out-Null
$d.Content.Find.Execute("74-48=", $true, $false, $false, $false, $false, $true, 1, $false, "30+16=", 2) | Out-Null
$d.Content.Find.Execute("35+34=", $true, $false, $false, $false, $false, $true, 1, $false, "82-80=", 2) | Out-Null
$d.Content.Find.Execute("57-26=", $true, $false, $false, $false, $false, $true, 1, $false, "25+43=", 2) | Out-Null
$d.Content.Find.Execute("26+5=", $true, $false, $false, $false, $false, $true, 1, $false, "46-41=", 2) | Out-Null
$d.Content.Find.Execute("88+2=", $true, $false, $false, $false, $false, $true, 1, $false, "50-40=", 2) | Out-Null
$d.Content.Find.Execute("12+41=", $true, $false, $false, $false, $false, $true, 1, $false, "33-4=", 2) | Out-Null
$d.Content.Find.Execute("79-12=", $true, $false, $false, $false, $false, $true, 1, $false, "44+7=", 2) | Out-Null
$d.Content.Find.Execute("82-26=", $true, $false, $false, $false, $false, $true, 1, $false, "47-9=", 2) | Out-Null
$d.Content.Find.Execute("34-31=", $true, $false, $false, $false, $false, $true, 1, $false, "73+14=", 2) | Out-Null
$d.Content.Find.Execute("7+65=", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=", 2) | Out-Null
$d.Content.Find.Execute("49+4=", $true, $false, $false, $false, $false, $true, 1, $false, "79-29=", 2) | Out-Null
$d.Content.Find.Execute("10+63=", $true, $false, $false, $false, $false, $true, 1, $false, "17+66=", 2) | Out-Null
$d.Content.Find.Execute("27+52=", $true, $false, $false, $false, $false, $true, 1, $false, "75+1=", 2) | Out-Null
$d.Content.Find.Execute("20-17=", $true, $false, $false, $false, $false, $true, 1, $false, "70-5=", 2) | Out-Null
$d.Content.Find.Execute("25+70=", $true, $false, $false, $false, $false, $true, 1, $false, "19+1=", 2) | Out-Null
$d.Content.Find.Execute("61-6=", $true, $false, $false, $false, $false, $true, 1, $false, "20-8=", 2) | Out-Null
$d.Content.Find.Execute("18+2=", $true, $false, $false, $false, $false, $true, 1, $false, "29+21=", 2) | Out-Null
$d.Content.Find.Execute("37+55=", $true, $false, $false, $false, $false, $true, 1, $false, "96-85=", 2) | Out-Null
$d.Content.Find.Execute("52-10=", $true, $false, $false, $false, $false, $true, 1, $false, "88-1=", 2) | Out-Null
$d.Content.Find.Execute("72+5=", $true, $false, $false, $false, $false, $true, 1, $false, "52-3=", 2) | Out-Null
$d.Content.Find.Execute("56+0=", $true, $false, $false, $false, $false, $true, 1, $false, "43+51=", 2) | Out-Null
$d.Content.Find.Execute("51-31=", $true, $false, $false, $false, $false, $true, 1, $false, "95-15=", 2) | Out-Null
$d.Content.Find.Execute("85-60=", $true, $false, $false, $false, $false, $true, 1, $false, "99-40=", 2) | Out-Null
$d.Content.Find.Execute("70-58=", $true, $false, $false, $false, $false, $true, 1, $false, "85-43=", 2) | Out-Null
$d.Content.Find.Execute("5+37=", $true, $false, $false, $false, $false, $true, 1, $false, "29+19=", 2) | Out-Null
$d.Content.Find.Execute("89-25=", $true, $false, $false, $false, $false, $true, 1, $false, "17-12=", 2) | Out-Null
$d.Content.Find.Execute("38+30=", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=", 2) | Out-Null
$d.Content.Find.Execute("81-62=", $true, $false, $false, $false, $false, $true, 1, $false, "86+7=", 2) | Out-Null
$d.Content.Find.Execute("85-4=", $true, $false, $false, $false, $false, $true, 1, $false, "93-75=", 2) | Out-Null
$d.Content.Find.Execute("17+17=", $true, $false, $false, $false, $false, $true, 1, $false, "13+48=", 2) | Out-Null
$d.Content.Find.Execute("32-1=", $true, $false, $false, $false, $false, $true, 1, $false, "8+2=", 2) | Out-Null
$d.Content.Find.Execute("17-1=", $true, $false, $false, $false, $false, $true, 1, $false, "82-31=", 2) | Out-Null
$d.Content.Find.Execute("80-7=", $true, $false, $false, $false, $false, $true, 1, $false, "34+51=", 2) | Out-Null
$d.Content.Find.Execute("23-4=", $true, $false, $false, $false, $false, $true, 1, $false, "74-61=", 2) | Out-Null
$d.Content.Find.Execute("99-46=", $true, $false, $false, $false, $false, $true, 1, $false, "1+47=", 2) | Out-Null
$d.Content.Find.Execute("30+32=", $true, $false, $false, $false, $false, $true, 1, $false, "37+2=", 2) | Out-Null
$d.Content.Find.Execute("62+14=", $true, $false, $false, $false, $false, $true, 1, $false, "28-22=", 2) | Out-Null
$d.Content.Find.Execute("56-35=", $true, $false, $false, $false, $false, $true, 1, $false, "34+32=", 2) | Out-Null
